$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.528.43"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.470.00"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.77"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.90"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.850.34"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.00"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.456.21"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.560.42"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0942"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.11"
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.12"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.07"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.79"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.44"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.06"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.45"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.23"
$ws.Range("E35").Value = "  -4.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("E37").Value = "  -6.65%  "
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E40").Value = "  -4.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.945.69"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.00"
$ws.Range("E44").Value = "  -4.43%  "
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.10"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.708.79"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.39"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.18"
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.74"
$ws.Range("E51").Value = "  +3.47%  "
